$wb = $excel.ActiveWorkbook

# Sheet: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7472
$ws.Range("J3").Value = 7861
$ws.Range("C4").Value = 1845
$ws.Range("J4").Value = 1711
$ws.Range("J5").Value = 616
$ws.Range("J6").Value = 10736
$ws.Range("C7").Value = 28389
$ws.Range("J7").Value = 28396

# Sheet: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J7").Value = 811
$ws.Range("J8").Value = 1794
$ws.Range("J9").Value = 142
$ws.Range("J10").Value = 204
$ws.Range("J11").Value = 508
$ws.Range("J15").Value = 351
$ws.Range("J19").Value = 821
$ws.Range("J20").Value = 614
$ws.Range("J27").Value = 171
$ws.Range("J29").Value = 1512
$ws.Range("J30").Value = 99
$ws.Range("J33").Value = 1286
$ws.Range("J36").Value = 385
$ws.Range("J37").Value = 871
$ws.Range("J42").Value = 1207
$ws.Range("J47").Value = 206
$ws.Range("J52").Value = 723
$ws.Range("J55").Value = 446
$ws.Range("J60").Value = 167
$ws.Range("C63").Value = 274
$ws.Range("J63").Value = 87
$ws.Range("J65").Value = 714
$ws.Range("J67").Value = 1034
$ws.Range("J73").Value = 277
$ws.Range("J76").Value = 404
$ws.Range("J77").Value = 199
$ws.Range("J78").Value = 331
$ws.Range("J79").Value = 775
$ws.Range("J83").Value = 573
$ws.Range("J84").Value = 237
$ws.Range("J85").Value = 1170
$ws.Range("J88").Value = 302
$ws.Range("J90").Value = 298
$ws.Range("J91").Value = 325
$ws.Range("J94").Value = 316
$ws.Range("J95").Value = 406
$ws.Range("J96").Value = 322
$ws.Range("J97").Value = 256
$ws.Range("C101").Value = 28389
$ws.Range("J101").Value = 28396

# Sheet: West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 126
$ws.Range("J7").Value = 322

# Sheet: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J3").Value = 243
$ws.Range("J6").Value = 258
$ws.Range("J7").Value = 811

# Sheet: Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 141
$ws.Range("J7").Value = 508

# Sheet: South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 423
$ws.Range("J6").Value = 335
$ws.Range("J7").Value = 1170

# Sheet: Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 202
$ws.Range("J6").Value = 312
$ws.Range("J7").Value = 723

# Sheet: Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 517
$ws.Range("J6").Value = 667
$ws.Range("J7").Value = 1794

# Sheet: South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 171
$ws.Range("J7").Value = 573

# Sheet: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 292
$ws.Range("J7").Value = 1286

# Sheet: West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 144
$ws.Range("J7").Value = 406

# Sheet: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 264
$ws.Range("J3").Value = 291
$ws.Range("J7").Value = 871

# Sheet: New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 208
$ws.Range("J7").Value = 714

# Sheet: Fuller Park
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J2").Value = 36
$ws.Range("J7").Value = 99

# Sheet: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 391
$ws.Range("J7").Value = 1034

# Sheet: South Deering
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 237

# Sheet: Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 531
$ws.Range("J4").Value = 82
$ws.Range("J6").Value = 383
$ws.Range("J7").Value = 1512

# Sheet: Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 201
$ws.Range("J3").Value = 232
$ws.Range("J4").Value = 40
$ws.Range("J6").Value = 318
$ws.Range("J7").Value = 821

# Sheet: River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 73
$ws.Range("J6").Value = 208
$ws.Range("J7").Value = 404

# Sheet: Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J4").Value = 52
$ws.Range("J7").Value = 1207

# Sheet: Avondale
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 117
$ws.Range("J7").Value = 204

# Sheet: Rogers Park
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J2").Value = 88
$ws.Range("J3").Value = 100
$ws.Range("J7").Value = 331

# Sheet: Lower West Side
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 87
$ws.Range("J7").Value = 446

# Sheet: Washington Park
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 325

# Sheet: Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 222
$ws.Range("J7").Value = 775

# Sheet: Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 48
$ws.Range("J6").Value = 178
$ws.Range("J7").Value = 614

# Sheet: Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J6").Value = 114
$ws.Range("J7").Value = 385

# Sheet: West Loop
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 168
$ws.Range("J7").Value = 316

# Sheet: Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J3").Value = 53
$ws.Range("J7").Value = 206

# Sheet: Brighton Park
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 73
$ws.Range("J6").Value = 162
$ws.Range("J7").Value = 351

# Sheet: Avalon Park
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J5").Value = 6
$ws.Range("J7").Value = 142

# Sheet: Portage Park
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 73
$ws.Range("J7").Value = 277

# Sheet: West Town
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J6").Value = 177
$ws.Range("J7").Value = 256

# Sheet: United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 302

# Sheet: Edgewater
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 41
$ws.Range("J7").Value = 171

# Sheet: Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 81
$ws.Range("J7").Value = 298

# Sheet: Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J6").Value = 51
$ws.Range("J7").Value = 167

# Sheet: Riverdale
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 199
